# Insert a new data row above current row 66 (shifting rows 66-266 down to 67-267)
# and populate it with the new weekly record, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("66:66").Insert()

$ws.Range("A66").Value = 3
$ws.Range("B66").Value = "Femacal de La Calera"
$ws.Range("C66").Value = "Coquimbo"
$ws.Range("D66").Value = 44607
$ws.Range("E66").Value = 5
$ws.Range("F66").Value = 100112039
$ws.Range("G66").Value = "Ciboulette"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 120
$ws.Range("K66").Value = 1500
$ws.Range("L66").Value = 1500
$ws.Range("M66").Value = 1500
$ws.Range("N66").Value = "`$/docena de atados"
$ws.Range("O66").Value = "Provincia de Quillota"
$ws.Range("P66").Value = 500
$ws.Range("Q66").Value = 3
$ws.Range("R66").Value = "Hortaliza"
